# Quarterly database update: roll the 10-quarter window forward by one
# quarter (drop the oldest quarter "فصل چهارم منتهی به 1399/06", shift all
# remaining quarters left, and append the new quarter
# "فصل دوم منتهی به 1401/12" with freshly-reported figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (quarter labels) -----------------------------------
# Row 8 and Row 24 both show the same rolling 10-quarter header strip.
$quarters = @(
    "فصل اول منتهی به 1399/09",
    "فصل دوم منتهی به 1399/12",
    "فصل سوم منتهی به 1400/03",
    "فصل چهارم منتهی به 1400/06",
    "فصل اول منتهی به 1400/09",
    "فصل دوم منتهی به 1400/12",
    "فصل سوم منتهی به 1401/03",
    "فصل چهارم منتهی به 1401/06",
    "فصل اول منتهی به 1401/09",
    "فصل دوم منتهی به 1401/12"
)

$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $ws.Range("E$r").Value = $quarters[0]
    $ws.Range("F$r").Value = $quarters[1]
    $ws.Range("G$r").Value = $quarters[2]
    $ws.Range("H$r").Value = $quarters[3]
    $ws.Range("I$r").Value = $quarters[4]
    $ws.Range("J$r").Value = $quarters[5]
    $ws.Range("K$r").Value = $quarters[6]
    $ws.Range("L$r").Value = $quarters[7]
    $ws.Range("M$r").Value = $quarters[8]
    $ws.Range("N$r").Value = $quarters[9]
}

# --- Data rows --------------------------------------------------------
# Each row: new E..N values = old F..N values shifted left, with a new
# figure for the newly-added quarter appended at N.
$rowValues = @{
    10 = @(49208, 29466, 74866, 14969, 2340, 5638, 14936, 143726, 144936, 262694)
    16 = @(6656, 5977, 4698, 4615, 4469, 3865, 1952, 273, 1516, 4328)
    17 = @(48308, 57804, -13086, 170139, 79582, 79304, 87012, 125248, 129267, 136091)
    19 = @(3727051, 3732080, 3071725, -1392637, 2178988, 1252691, 1568423, 33298696, 2915396, 43263697)
    20 = @(3831223, 3825327, 3138203, -1202914, 2265379, 1341498, 1672323, 33567943, 3191115, 43666810)
    26 = @(143, 172, 183, 199, 178, 234, 234, 240, 245, 248)
    27 = @(549, 520, 520, 500, 514, 554, 554, 560, 561, 560)
}

$dataCols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range("$($dataCols[$i])$r").Value = $vals[$i]
    }
}
